$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.607.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.868.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.37%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.53%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'338.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.96%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.63%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4685"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +4.82%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'47.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08034"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.93%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.11%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.047"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.859.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.289"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.61%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.96%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001042"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.58%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06611"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.30%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.62%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'28.629.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.36%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.488"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.20%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.29%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.252"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.083.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.83%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'160.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.73%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +2.83%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.487"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.83%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'119.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9724"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09526"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.59%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.592"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.382"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.72%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.363"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.06203"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.85%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02254"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.48%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.78%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.182"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5934"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.9984"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.68%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1881"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.80%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'10.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.27%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.05%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.58%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Decentraland"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.5557"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.98%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.07423"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +12.27%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.956"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.04%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.087"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +13.89%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'112.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.92%  "
$ws.Range("E51").Style = "Normal"
